$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows above row 444, shifting the existing rows 444-451 down to 448-455.
$ws.Range("A444:A447").EntireRow.Insert(-4121) | Out-Null

# Row 444 (new)
$ws.Range("A444").Value = 10
$ws.Range("B444").Value = "Vega Modelo de Temuco"
$ws.Range("C444").Value = "La Araucanía"
$ws.Range("D444").Value = 44595
$ws.Range("E444").Value = 9
$ws.Range("F444").Value = "Fruta"
$ws.Range("G444").Value = 100103
$ws.Range("H444").Value = "Frutos de hueso (carozo)"
$ws.Range("I444").Value = 100103006
$ws.Range("J444").Value = "Nectarín"
$ws.Range("K444").Value = "June Pearl"
$ws.Range("L444").Value = "Primera"
$ws.Range("M444").Value = 5
$ws.Range("N444").Value = 410000
$ws.Range("O444").Value = 410000
$ws.Range("P444").Value = 410000
$ws.Range("Q444").Value = "$/bins (420 kilos)"
$ws.Range("R444").Value = "Región de O'Higgins"
$ws.Range("S444").Value = 976
$ws.Range("T444").Value = 420

# Row 445 (new)
$ws.Range("A445").Value = 10
$ws.Range("B445").Value = "Vega Modelo de Temuco"
$ws.Range("C445").Value = "La Araucanía"
$ws.Range("D445").Value = 44595
$ws.Range("E445").Value = 9
$ws.Range("F445").Value = "Fruta"
$ws.Range("G445").Value = 100103
$ws.Range("H445").Value = "Frutos de hueso (carozo)"
$ws.Range("I445").Value = 100103006
$ws.Range("J445").Value = "Nectarín"
$ws.Range("K445").Value = "Ruby Diamond"
$ws.Range("L445").Value = "Primera"
$ws.Range("M445").Value = 255
$ws.Range("N445").Value = 15000
$ws.Range("O445").Value = 16000
$ws.Range("P445").Value = 15490
$ws.Range("Q445").Value = "$/bandeja 18 kilos granel"
$ws.Range("R445").Value = "Región de O'Higgins"
$ws.Range("S445").Value = 861
$ws.Range("T445").Value = 18

# Row 446 (new)
$ws.Range("A446").Value = 10
$ws.Range("B446").Value = "Vega Modelo de Temuco"
$ws.Range("C446").Value = "La Araucanía"
$ws.Range("D446").Value = 44595
$ws.Range("E446").Value = 9
$ws.Range("F446").Value = "Fruta"
$ws.Range("G446").Value = 100103
$ws.Range("H446").Value = "Frutos de hueso (carozo)"
$ws.Range("I446").Value = 100103006
$ws.Range("J446").Value = "Nectarín"
$ws.Range("K446").Value = "Super Queen"
$ws.Range("L446").Value = "Primera"
$ws.Range("M446").Value = 190
$ws.Range("N446").Value = 16000
$ws.Range("O446").Value = 18000
$ws.Range("P446").Value = 16684
$ws.Range("Q446").Value = "$/bandeja 18 kilos granel"
$ws.Range("R446").Value = "Región de O'Higgins"
$ws.Range("S446").Value = 927
$ws.Range("T446").Value = 18

# Row 447 (new)
$ws.Range("A447").Value = 10
$ws.Range("B447").Value = "Vega Modelo de Temuco"
$ws.Range("C447").Value = "La Araucanía"
$ws.Range("D447").Value = 44595
$ws.Range("E447").Value = 9
$ws.Range("F447").Value = "Fruta"
$ws.Range("G447").Value = 100103
$ws.Range("H447").Value = "Frutos de hueso (carozo)"
$ws.Range("I447").Value = 100103006
$ws.Range("J447").Value = "Nectarín"
$ws.Range("K447").Value = "Venus"
$ws.Range("L447").Value = "Primera"
$ws.Range("M447").Value = 8
$ws.Range("N447").Value = 390000
$ws.Range("O447").Value = 390000
$ws.Range("P447").Value = 390000
$ws.Range("Q447").Value = "$/bins (420 kilos)"
$ws.Range("R447").Value = "Región de O'Higgins"
$ws.Range("S447").Value = 929
$ws.Range("T447").Value = 420
